$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two API endpoint URLs (moving from local dev URLs to deployed
# Vercel URLs) in column C.
$ws.Range("C36").Value = "https://tds-ga3-7.vercel.app/similarity"
$ws.Range("C41").Value = "https://tds-ga4-3.vercel.app/api/outline"

# Move the selection / scroll position to match where the author ended up
# (row 41, with the view scrolled down to around row 30).
$ws.Range("C41").Select()
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 3
